# Rename the two worksheets to reflect the unified DataNode / DataTable / Entity
# naming convention, and switch the active (selected) tab from the first sheet
# to the second sheet, matching the author's re-save of the workbook.

$wb = $excel.ActiveWorkbook

$wsDataNode  = $wb.Worksheets.Item(1)   # was "Property1"
$wsDataTable = $wb.Worksheets.Item(2)   # was "Record"

$wsDataNode.Name  = "DataNode"
$wsDataTable.Name = "DataTable"

# Make the second sheet ("DataTable") the active / selected tab.
$wsDataTable.Activate()
